$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracker rows 135-139 (linear-interpolation fill for missing daily
# updates where the underlying report data is static), mirroring the
# existing row layout: A=tracker_date, B=report_date, C..L=counts,
# M=source url.

$rows = @(
    @{ Row=135; A="19.02.2024"; B="16.02.2024"; C=28775; D=12300; E=8400; F=68552; G=8663; H=6327; I=7000; J=395; K=105; L=4450; M="https://web.archive.org/web/20240219021050/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=136; A="20.02.2024"; B="16.02.2024"; C=28775; D=12300; E=8400; F=68552; G=8663; H=6327; I=7000; J=395; K=105; L=4450; M="https://web.archive.org/web/20240220031301/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=137; A="21.02.2024"; B="21.02.224";  C=29313; D=12300; E=8400; F=69333; G=8663; H=6327; I=7000; J=395; K=105; L=4450; M="https://web.archive.org/web/20240221120805/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=138; A="22.02.2024"; B="21.02.224";  C=29313; D=12300; E=8400; F=69333; G=8663; H=6327; I=7000; J=395; K=105; L=4450; M="https://web.archive.org/web/20240222212239/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=139; A="23.02.2024"; B="21.02.224";  C=29313; D=12300; E=8400; F=69333; G=8663; H=6327; I=7000; J=395; K=105; L=4450; M="https://web.archive.org/web/20240223181530/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}

# Keep the selection / active cell on the newly appended last row, matching
# the tracker's append-only editing pattern.
$ws.Range("M139").Select() | Out-Null
